$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -13.127
$ws.Range("C4").Value = -14.55030000000001
$ws.Range("E6").Value = 12.27649999999999
$ws.Range("C7").Value = -11.6825
$ws.Range("E7").Value = 12.8177
$ws.Range("C8").Value = -12.00909999999999
$ws.Range("E8").Value = 13.38160000000001
$ws.Range("A11").Value = -21.88100000000003
$ws.Range("A12").Value = -21.43930000000001
$ws.Range("C12").Value = -11.8622
$ws.Range("C14").Value = -11.47999999999999
$ws.Range("A15").Value = -21.26790000000001
$ws.Range("E19").Value = 12.76289999999999
$ws.Range("E21").Value = 12.87149999999999
$ws.Range("C22").Value = -11.05549999999998
$ws.Range("E24").Value = 12.79969999999998
$ws.Range("E25").Value = 13.40530000000001
